# Insert a new weekly price record as the first data row for this
# "Hortaliza, Vega Central Mapocho de Santiago - Orégano" sheet.
#
# The sheet has a header row (row 1) followed by data rows starting at
# row 2 (sorted with the most-recent weekly entries near the top of the
# block that starts at row 55). This edit inserts one brand-new record
# at row 55, pushing the previous rows 55-116 down to 56-117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 55 (shifts rows 55:116 down to 56:117)
$ws.Rows.Item(55).Insert()

# Populate the new row 55 with the new weekly record
$ws.Cells.Item(55, 1).Value  = 9
$ws.Cells.Item(55, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(55, 3).Value  = "Metropolitana"
$ws.Cells.Item(55, 4).Value  = 45175
$ws.Cells.Item(55, 5).Value  = 13
$ws.Cells.Item(55, 6).Value  = 100112029
$ws.Cells.Item(55, 7).Value  = "Orégano"
$ws.Cells.Item(55, 8).Value  = "Sin especificar"
$ws.Cells.Item(55, 9).Value  = "Primera"
$ws.Cells.Item(55, 10).Value = 16
$ws.Cells.Item(55, 11).Value = 21000
$ws.Cells.Item(55, 12).Value = 21000
$ws.Cells.Item(55, 13).Value = 21000
$ws.Cells.Item(55, 14).Value = "$/docena de atados"
$ws.Cells.Item(55, 15).Value = "Región Metropolitana"
$ws.Cells.Item(55, 16).Value = 7000
$ws.Cells.Item(55, 17).Value = 3
$ws.Cells.Item(55, 18).Value = "Hortaliza"

# Keep the date column formatted like the rest of column D
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(56, 4).NumberFormat

Write-Host "Inserted new row 55; sheet now spans" $ws.UsedRange.Address()
